# Update column F (dSF) values on the active worksheet to reflect
# repulled data / recalculated means, per the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -1
    3  = 2
    4  = 1
    5  = -1
    6  = 1
    7  = 2
    8  = 3
    9  = -4
    10 = 2
    11 = 0
    12 = 5
    13 = 2
    14 = 1
    15 = 1
    16 = -5
    17 = 4
    18 = -4
    19 = 4
    20 = -3
    21 = -2
    22 = 3
    23 = 5
    24 = 6
    25 = -9
    26 = -4
    27 = -5
    28 = -3
    30 = -1
    31 = 4
    32 = 1
    33 = -3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
